$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string
$ws.Range("A1").Value = "Datos actualizados a 14 de Julio de 2020 a las 07:21"

# Row 15 (Pakistan)
$ws.Range("B15").Value = 253604
$ws.Range("C15").Value = 1979
$ws.Range("D15").Value = 170656
$ws.Range("E15").Value = 77628
$ws.Range("G15").Value = 54
$ws.Range("H15").Value = 5320

# Row 72 (Kirguistan)
$ws.Range("B72").Value = 11444
$ws.Range("C72").Value = 327
$ws.Range("D72").Value = 3538
$ws.Range("E72").Value = 7757
$ws.Range("G72").Value = 2
$ws.Range("H72").Value = 149

# Row 103 (Tailandia)
$ws.Range("B103").Value = 3227
$ws.Range("C103").Value = 7
$ws.Range("D103").Value = 3091
$ws.Range("E103").Value = 78

# Row 186 (Butan)
$ws.Range("D186").Value = 78
$ws.Range("E186").Value = 6
